# Normalize "Recorded By" (column G) entries so that a leading "System"
# entry in the comma-separated list is moved to the end of the list.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#      "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text -split ", "
    if ($parts.Length -gt 1 -and $parts[0] -eq "System") {
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + @("System")
        $newText = [string]::Join(", ", $newParts)
        $cell.Value2 = $newText
    }
}
